$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new log rows (125 -> 126, 126 -> 127) to the feed_logs sheet.
$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 2).Value = 1
$ws.Cells.Item(126, 3).Value = "2024-06-17 14:12:32"
$ws.Cells.Item(126, 4).Value = 200
$ws.Cells.Item(126, 5).Value = 23

$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 2).Value = 2
$ws.Cells.Item(127, 3).Value = "2024-06-17 14:12:33"
$ws.Cells.Item(127, 4).Value = 200
$ws.Cells.Item(127, 5).Value = 1
